$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "58.032.32"
$ws.Range("E2").Value = "  +2.37%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.067.27"
$ws.Range("E3").Value = "  +1.60%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "517.01"
$ws.Range("E5").Value = "  +1.28%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "141.84"
$ws.Range("E6").Value = "  +1.78%  "
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("E8").Value = "  +1.09%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "7.31"
$ws.Range("E9").Value = "  +2.77%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.108"
$ws.Range("E10").Value = "  +0.06%  "
$ws.Range("E11").Value = "  +2.00%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "3.592.78"
$ws.Range("E12").Value = "  +1.47%  "
$ws.Range("E13").Value = "  +2.86%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "26.39"
$ws.Range("E14").Value = "  +4.59%  "
$ws.Range("E15").Value = "  +1.22%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "58.058.45"
$ws.Range("E16").Value = "  +2.45%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.066.24"
$ws.Range("E17").Value = "  +1.58%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.11"
$ws.Range("E18").Value = "  +2.95%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.84"
$ws.Range("E19").Value = "  -2.06%  "
$ws.Range("E20").Value = "  +1.71%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "331.84"
$ws.Range("E21").Value = "  -0.45%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.00"
$ws.Range("E22").Value = "  -0.13%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.502"
$ws.Range("E23").Value = "  +0.09%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "65.42"
$ws.Range("E24").Value = "  +0.84%  "
$ws.Range("E25").Value = "  +2.43%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.999"
$ws.Range("E26").Value = "  -0.14%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0₃0909"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.48"
$ws.Range("E28").Value = "  +1.98%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.24"
$ws.Range("E29").Value = "  +6.32%  "
$ws.Range("E30").Value = "  +1.05%  "
$ws.Range("E31").Value = "  +4.32%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "20.65"
$ws.Range("E32").Value = "  +1.14%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "155.18"
$ws.Range("B34").Value = "EnergySwap"
$ws.Range("C34").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "27.63"
$ws.Range("E34").Value = "  +1.77%  "
$ws.Range("B35").Value = "NEARProtocol"
$ws.Range("C35").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.53"
$ws.Range("E35").Value = "  +1.32%  "
$ws.Range("E36").Value = "  +3.40%  "
$ws.Range("E37").Value = "  +3.73%  "
$ws.Range("E38").Value = "  +2.63%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.109.54"
$ws.Range("E39").Value = "  +1.78%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.92"
$ws.Range("E40").Value = "  +2.47%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "36.79"
$ws.Range("E41").Value = "  +1.20%  "
$ws.Range("E42").Value = "  +0.01%  "
$ws.Range("E43").Value = "  +0.30%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.301.99"
$ws.Range("E44").Value = "  +4.59%  "
$ws.Range("E45").Value = "  +4.24%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.39"
$ws.Range("E46").Value = "  +2.53%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "20.89"
$ws.Range("E47").Value = "  +5.74%  "
$ws.Range("E48").Value = "  +1.42%  "
$ws.Range("E49").Value = "  +1.77%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.740"
$ws.Range("E50").Value = "  +9.95%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "253.90"
$ws.Range("E51").Value = "  +9.87%  "
